# Applies the "Last version" edit to data_vuv.xlsx:
#  - H61: 60075 -> -999
#  - Append new row 71 with a new UFO observation entry
#  - Update the active selection to C44 (closest reproducible approximation
#    of the recorded view state change; this COM runtime does not expose a
#    way to persist window scroll position / frozen panes / zoom level)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- existing-cell edit -----------------------------------------------
$ws.Range("H61").Value = -999

# --- new row 71 ----------------------------------------------------------
$ws.Range("A71").Value = 379
$ws.Range("B71").Value = 97851
$ws.Range("C71").Value = "KLDT-O5WB"

$ws.Range("D71").Value = 49225
$ws.Range("D71").NumberFormat = "#,##0"

$ws.Range("E71").Value = "KLDT-E5WD"
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = "Bright flash at the top of the outboard side associated with a titanium signal in KT2."
$ws.Range("H71").Value = 747
$ws.Range("I71").Value = "There's yet another UFO right after, not as visible."
# J71 intentionally left blank (no Initial_pos recorded for this entry)
$ws.Range("K71").Value = "Ti"
$ws.Range("L71").Value = "Based on comment (KT2)"
$ws.Range("M71").Value = "4D"

# --- selection / view ------------------------------------------------
$ws.Range("C44").Select()
